# "Update in the simulation" — trims Sheet1's dataset down to the 5
# non-hidden simulation rows (dropping the hidden detail rows that held
# the old 1000-transaction sweep), refreshes the still-visible rows with
# the new 95000-transaction parameters, un-hides column E, re-points the
# Sheet1 charts at the smaller range, and makes Sheet1 the active tab.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ------------------------------------------------------------------
# 1. Drop the hidden "detail" rows from Sheet1. Row 10 (not hidden,
#    "SDT 200X200") slides up into row 6 once rows 6:9 go away; the
#    remaining hidden rows 11:20 (now 7:16) are removed next.
# ------------------------------------------------------------------
$ws1.Rows("6:9").Delete()
$ws1.Rows("7:16").Delete()

# ------------------------------------------------------------------
# 2. The simulation's transaction-count parameter moved from 1000 to
#    95000 for the two runs that remain visible.
# ------------------------------------------------------------------
$ws1.Range("E5").Value = "[100,100,95000,0.0015]"
$ws1.Range("E6").Value = "[200,200,95000,0.0015]"

# ------------------------------------------------------------------
# 3. Column E ("Param") is no longer hidden.
# ------------------------------------------------------------------
$ws1.Columns("E").Hidden = $false

# ------------------------------------------------------------------
# 4. Re-point the three Sheet1 charts at the shrunk table (rows 2-6
#    instead of 2-20) while keeping their series names/order intact.
# ------------------------------------------------------------------
$chart1 = $ws1.ChartObjects().Item(1).Chart
$chart1.SeriesCollection().Item(1).Formula = '=SERIES("Block Time",Sheet1!$B$2:$B$6,Sheet1!$C$2:$C$6,1)'
$chart1.SeriesCollection().Item(2).Formula = '=SERIES("Fee Collected",Sheet1!$B$2:$B$6,Sheet1!$F$2:$F$6,2)'

$chart2 = $ws1.ChartObjects().Item(2).Chart
$chart2.SeriesCollection().Item(1).Formula = '=SERIES(Sheet1!$F$1,Sheet1!$B$2:$B$6,Sheet1!$F$2:$F$6,1)'

$chart3 = $ws1.ChartObjects().Item(3).Chart
$chart3.SeriesCollection().Item(1).Formula = '=SERIES(Sheet1!$C$1,Sheet1!$B$2:$B$6,Sheet1!$C$2:$C$6,1)'
$chart3.SeriesCollection().Item(2).Formula = '=SERIES(Sheet1!$D$1,Sheet1!$B$2:$B$6,Sheet1!$D$2:$D$6,2)'

# The deleted rows were all zero-height (hidden), so the charts' actual
# on-sheet position/size doesn't change — but this engine doesn't shift
# the drawing anchors along with the row delete, so pin them back to
# the (unchanged) absolute Top/Height the anchors resolved to before.
$co1 = $ws1.ChartObjects().Item(1)
$co2 = $ws1.ChartObjects().Item(2)
$co3 = $ws1.ChartObjects().Item(3)
$co1.Height = 311.811023622047
$co2.Height = 311.811023622047
$co3.Top    = 327.5
$co3.Height = 311.811023622047

# ------------------------------------------------------------------
# 5. Sheet1 becomes the active sheet/tab (it was Sheet1 (2) before),
#    with the selection left on E7, just past the trimmed table.
# ------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("E7").Select()
